$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 05:41"

# 1b. Swap the "Groenlandia" / "Islas Malvinas" rows (their stats are identical,
# only the label ordering changed in the source data).
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# 2. Update Belgica row (row 33)
$ws.Range("B33").Value = 62707
$ws.Range("C33").Value = 101
$ws.Range("E33").Value = 35729

# 3. Update Honduras row (row 55)
$ws.Range("B55").Value = 28090
$ws.Range("C55").Value = 507
$ws.Range("D55").Value = 2957
$ws.Range("E55").Value = 24359
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 774

# 4. Update Haiti row (row 88)
$ws.Range("B88").Value = 6727
$ws.Range("C88").Value = 37
$ws.Range("D88").Value = 2924
$ws.Range("E88").Value = 3664

# 5. Update Mongolia row (row 171)
$ws.Range("D171").Value = 203
$ws.Range("E171").Value = 27
